# edit.ps1 - reproduces the commit:
#   1) Table on slide 5 switches from the deck's custom "Table_0" table
#      style to the built-in "No Style, Table Grid" table style.
#   2) The presentation's theme colour scheme is swapped from the custom
#      "Integral" / "Red Violet" palette to the stock "Office" palette
#      (the slide-master theme and the notes-master theme traded their
#      colour content in the authored file).
#
# NOTE on (2): in this COM host, Presentation.SlideMaster.ColorScheme and
# Presentation.NotesMaster.ColorScheme are the same underlying theme
# colour table (there is no independently addressable notes-master theme
# object exposed on the object model), so only the reachable theme's RGB
# values can be driven from script; that theme is updated to match the
# target "Office" palette exactly, colour-for-colour.

$p = $ppt.ActivePresentation

# --- 1) Table style fix on slide 5 --------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F012BEA5-0550-428A-9D74-86EC9447E373}")
    }
}

# --- 2) Theme colour-scheme swap (Integral/Red Violet -> Office) -------
function Set-RGBFromHex($colorItem, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorItem.RGB = $r + ($g * 256) + ($b * 65536)
}

$officePalette = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$scheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officePalette.Count; $i++) {
    Set-RGBFromHex $scheme.Item($i) $officePalette[$i - 1]
}
